$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric cells
$ws.Range("A3").Value = 112564182
$ws.Range("B3").Value = 103742
$ws.Range("E3").Value = 340
$ws.Range("Q3").Value = 542782
$ws.Range("R3").Value = 6404993
$ws.Range("S3").Value = 50

# Plain text cells
$ws.Range("C3").Value = "Ovaliderad"
$ws.Range("D3").Value = "EN"
$ws.Range("F3").Value = "Ryl"
$ws.Range("G3").Value = "Chimaphila umbellata"
$ws.Range("H3").Value = "(L.) W. P. C. Barton"
$ws.Range("J3").Value = "stj" + [char]0x00E4 + "lkar/str" + [char]0x00E5 + "n/skott"
$ws.Range("P3").Value = "N" + [char]0x00E4 + "set, 450 m SSV om, Sm"
$ws.Range("T3").Value = "Kalmar"
$ws.Range("U3").Value = "Vimmerby"
$ws.Range("V3").Value = "Sm" + [char]0x00E5 + "land"
$ws.Range("W3").Value = "S" + [char]0x00F6 + "dra Vi"
$ws.Range("X3").Value = "Hf-Vim-1025"
$ws.Range("AC3").Value = "S" + [char]0x00F6 + "ren Mj" + [char]0x00F6 + "sberg"
$ws.Range("AI3").Value = "Tallskog"
$ws.Range("AW3").Value = [char]0x00C5 + "ke R" + [char]0x00FC + "hling"
$ws.Range("AX3").Value = "Via " + [char]0x00C5 + "ke R" + [char]0x00FC + "hling"
$ws.Range("AY3").Value = "Florav" + [char]0x00E4 + "kteri Sverige"

# Text cells that must stay text (numeric-looking / date-looking strings)
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "20"

$ws.Range("Y3").NumberFormat = "@"
$ws.Range("Y3").Value = "2023-07-19"

$ws.Range("AA3").NumberFormat = "@"
$ws.Range("AA3").Value = "2023-07-19"

# Empty text placeholder cell (kept present but empty, like the source row)
$ws.Range("AT3").NumberFormat = "@"
$ws.Range("AT3").Value = ""

# Boolean cells
$ws.Range("AD3").Value = $False
$ws.Range("AE3").Value = $False
$ws.Range("AG3").Value = $False
